# Weekly price update: insert two new rows with the latest week's data
# (pushing the prior week's two rows down), then update the now-current
# rows 105 and 106 with the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 105 so the existing rows 105-106 become 107-108.
$ws.Rows.Item(105).Resize(2).Insert()

# Row 105 : Americana (o) - updated figures for the newer date (44706)
$ws.Cells.Item(105, 1).Value = 11
$ws.Cells.Item(105, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(105, 3).Value = "Bíobío"
$ws.Cells.Item(105, 4).Value = 44706
$ws.Cells.Item(105, 5).Value = 8
$ws.Cells.Item(105, 6).Value = 100112021
$ws.Cells.Item(105, 7).Value = "Ají"
$ws.Cells.Item(105, 8).Value = "Americana (o)"
$ws.Cells.Item(105, 9).Value = "Primera"
$ws.Cells.Item(105, 10).Value = 100
$ws.Cells.Item(105, 11).Value = 35000
$ws.Cells.Item(105, 12).Value = 36000
$ws.Cells.Item(105, 13).Value = 35500
$ws.Cells.Item(105, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(105, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(105, 16).Value = 1420
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"

# Row 106 : Chilena(o) - updated figures for the newer date (44706)
$ws.Cells.Item(106, 1).Value = 11
$ws.Cells.Item(106, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(106, 3).Value = "Bíobío"
$ws.Cells.Item(106, 4).Value = 44706
$ws.Cells.Item(106, 5).Value = 8
$ws.Cells.Item(106, 6).Value = 100112021
$ws.Cells.Item(106, 7).Value = "Ají"
$ws.Cells.Item(106, 8).Value = "Chilena(o)"
$ws.Cells.Item(106, 9).Value = "Primera"
$ws.Cells.Item(106, 10).Value = 40
$ws.Cells.Item(106, 11).Value = 45000
$ws.Cells.Item(106, 12).Value = 46000
$ws.Cells.Item(106, 13).Value = 45500
$ws.Cells.Item(106, 14).Value = "$/caja 25 kilos"
$ws.Cells.Item(106, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(106, 16).Value = 1820
$ws.Cells.Item(106, 17).Value = 25
$ws.Cells.Item(106, 18).Value = "Hortaliza"
